# "first draft of finished episode"
#
# Slide 4 (the communication/packed-memory diagram) has a second row of
# small rectangles representing x[]/y[] array cells. Several of those
# rectangles were purely decorative "blank" placeholders (hatch-pattern
# fill, no text) that are being removed; the remaining labeled cells in
# that row are shifted left/up to close up the resulting gaps so the row
# reads as a contiguous strip again.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)

function Get-ShapeById($slide, $id) {
    foreach ($shp in $slide.Shapes) {
        if ($shp.Id -eq $id) { return $shp }
    }
    return $null
}

# Shape.Left/Top round-trip through a single-precision (float32) COM
# property, so naively assigning target_emu/12700.0 can land one EMU off
# after PowerPoint re-quantizes it. Nudge the assigned point value within
# a tiny neighborhood until the value that comes back converts to exactly
# the EMU we want.
function Set-ExactEmu($shp, $axis, $targetEmu) {
    $targetPts = $targetEmu / 12700.0

    if ($axis -eq "Left") { $shp.Left = $targetPts } else { $shp.Top = $targetPts }
    $back = if ($axis -eq "Left") { $shp.Left } else { $shp.Top }
    if ([Math]::Round($back * 12700.0) -eq $targetEmu) { return }

    for ($i = 1; $i -le 4000; $i++) {
        foreach ($sign in @(1, -1)) {
            $cand = $targetPts + ($sign * $i * 0.0000001)
            if ($axis -eq "Left") { $shp.Left = $cand } else { $shp.Top = $cand }
            $back2 = if ($axis -eq "Left") { $shp.Left } else { $shp.Top }
            if ([Math]::Round($back2 * 12700.0) -eq $targetEmu) { return }
        }
    }
    Write-Host "WARN: could not hit exact EMU $targetEmu on axis $axis"
}

function Move-ShapeTo($id, $x, $y) {
    $shp = Get-ShapeById $s $id
    Set-ExactEmu $shp "Left" $x
    Set-ExactEmu $shp "Top"  $y
}

# --- Reposition the surviving second-row rectangles ---

Move-ShapeTo 37 2508316 659294   # x[1]
Move-ShapeTo 38 2830467 659293   # x[2]
Move-ShapeTo 40 4436571 659293   # y[2]
Move-ShapeTo 42 4114420 659293   # y[1]
Move-ShapeTo 43 3152618 659293   # x[3]
Move-ShapeTo 45 4758722 659293   # blank (unlabeled, solid fill)
Move-ShapeTo 4  5080873 659291   # y[4]
Move-ShapeTo 5  1867112 659295   # "4"
Move-ShapeTo 7  3470118 659293   # "5"
Move-ShapeTo 36 2186165 659295   # x[0]
Move-ShapeTo 41 3792269 659294   # y[0]

# --- Remove the decorative hatch-pattern placeholder rectangles ---

(Get-ShapeById $s 44).Delete()
(Get-ShapeById $s 8).Delete()
(Get-ShapeById $s 9).Delete()
(Get-ShapeById $s 39).Delete()

# --- Fix up shape stacking order to match the new arrangement ---

# "x[3]" (id 43) ends up re-inserted at the very end of the shape tree.
$shp43 = Get-ShapeById $s 43
for ($k = 0; $k -lt 9; $k++) {
    $shp43.ZOrder(2)   # msoBringForward
}

# "y[4]" duplicate (id 3) moves up, right after the blank rectangle (id 45)
# and before the two "Packed"/"Memory" text boxes.
$shp3 = Get-ShapeById $s 3
for ($k = 0; $k -lt 2; $k++) {
    $shp3.ZOrder(3)    # msoSendBackward
}
